$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$style_D2 = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "71.633.79"
$ws.Range("D2").Style = $style_D2
$ws.Range("E2").Value = "  +2.78%  "
$style_D3 = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "4.030.75"
$ws.Range("D3").Style = $style_D3
$ws.Range("E3").Value = "  +2.30%  "
$ws.Range("E4").Value = "  +0.00%  "
$style_D5 = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "525.67"
$ws.Range("D5").Style = $style_D5
$ws.Range("E5").Value = "  -0.98%  "
$style_D6 = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.56"
$ws.Range("D6").Style = $style_D6
$ws.Range("E6").Value = "  +1.39%  "
$style_D7 = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.626"
$ws.Range("D7").Style = $style_D7
$ws.Range("E7").Value = "  +0.62%  "
$ws.Range("E8").Value = "  +0.22%  "
$style_D9 = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.740"
$ws.Range("D9").Style = $style_D9
$ws.Range("E9").Value = "  +1.26%  "
$ws.Range("E10").Value = "  +1.43%  "
$style_D11 = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000342"
$ws.Range("D11").Style = $style_D11
$ws.Range("E11").Value = "  -0.53%  "
$style_D12 = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "46.26"
$ws.Range("D12").Style = $style_D12
$ws.Range("E12").Value = "  +7.94%  "
$style_D13 = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.83"
$ws.Range("D13").Style = $style_D13
$ws.Range("E13").Value = "  +2.92%  "
$style_D14 = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.670.00"
$ws.Range("D14").Style = $style_D14
$ws.Range("E14").Value = "  +2.11%  "
$style_D15 = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.036.56"
$ws.Range("D15").Style = $style_D15
$ws.Range("E15").Value = "  +2.31%  "
$style_D16 = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.43"
$ws.Range("D16").Style = $style_D16
$ws.Range("E16").Value = "  +7.75%  "
$style_D17 = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.33"
$ws.Range("D17").Style = $style_D17
$ws.Range("E17").Value = "  +1.04%  "
$ws.Range("E18").Value = "  +1.48%  "
$ws.Range("E19").Value = "  -1.60%  "
$style_D20 = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.638.17"
$ws.Range("D20").Style = $style_D20
$ws.Range("E20").Value = "  +2.83%  "
$style_D21 = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "443.18"
$ws.Range("D21").Style = $style_D21
$ws.Range("E21").Value = "  +1.79%  "
$style_D22 = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.59"
$ws.Range("D22").Style = $style_D22
$ws.Range("E22").Value = "  +5.23%  "
$style_D23 = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "95.03"
$ws.Range("D23").Style = $style_D23
$ws.Range("E23").Value = "  +7.41%  "
$style_D24 = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.41"
$ws.Range("D24").Style = $style_D24
$ws.Range("E24").Value = "  -1.08%  "
$style_D25 = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.26"
$ws.Range("D25").Style = $style_D25
$ws.Range("E25").Value = "  +2.67%  "
$ws.Range("E26").Value = "  -0.39%  "
$style_D27 = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.11"
$ws.Range("D27").Style = $style_D27
$style_D28 = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "37.20"
$ws.Range("D28").Style = $style_D28
$ws.Range("E28").Value = "  +1.13%  "
$style_D29 = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "13.70"
$ws.Range("D29").Style = $style_D29
$ws.Range("E29").Value = "  +2.38%  "
$style_D30 = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "699.31"
$ws.Range("D30").Style = $style_D30
$ws.Range("E30").Value = "  -1.20%  "
$ws.Range("E31").Value = "  +2.88%  "
$ws.Range("E32").Value = "  +1.78%  "
$style_D33 = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.94"
$ws.Range("D33").Style = $style_D33
$ws.Range("E33").Value = "  +14.37%  "
$style_D34 = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "68.01"
$ws.Range("D34").Style = $style_D34
$ws.Range("E34").Value = "  -0.26%  "
$style_D35 = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0904"
$ws.Range("D35").Style = $style_D35
$ws.Range("E35").Value = "  +3.31%  "
$ws.Range("E36").Value = "  +0.93%  "
$style_D37 = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "41.43"
$ws.Range("D37").Style = $style_D37
$ws.Range("E37").Value = "  +2.28%  "
$ws.Range("E38").Value = "  +3.79%  "
$style_D39 = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.54"
$ws.Range("D39").Style = $style_D39
$ws.Range("E39").Value = "  +17.67%  "
$ws.Range("E40").Value = "  +0.07%  "
$style_D41 = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0493"
$ws.Range("D41").Style = $style_D41
$ws.Range("E41").Value = "  +1.71%  "
$ws.Range("E42").Value = "  -0.10%  "
$style_D43 = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.86"
$ws.Range("D43").Style = $style_D43
$ws.Range("E43").Value = "  +0.68%  "
$ws.Range("E44").Value = "  +1.07%  "
$ws.Range("E45").Value = "  +3.80%  "
$ws.Range("E46").Value = "  +2.88%  "
$style_D47 = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.18"
$ws.Range("D47").Style = $style_D47
$ws.Range("E47").Value = "  -1.25%  "
$style_D48 = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.29"
$ws.Range("D48").Style = $style_D48
$ws.Range("E48").Value = "  +7.06%  "
$style_D49 = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.000278"
$ws.Range("D49").Style = $style_D49
$ws.Range("E49").Value = "  +17.48%  "
$style_D50 = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.39"
$ws.Range("D50").Style = $style_D50
$ws.Range("E50").Value = "  +0.84%  "
$ws.Range("E51").Value = "  -5.87%  "
